$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the indicator "Implications" text for the two MRIP recreational rows to
# remove references to preliminary 2024 data that excludes Nov/Dec.
$ws.Range("C4").Value = "Recent trip numbers are near an all-time high, but have decreased from 2023. Catch (not shown) generally reflects trip patterns. High regulatory complexity is likely contributing to recreational fishing trends."
$ws.Range("C5").Value = "The recreational black sea bass fishery has a catch-and-release component, and management measures are being implemented to reduce recreational harvest. "

# The rows are now shorter, so update the (auto-calculated) row heights to
# reflect the reduced amount of wrapped text.
$ws.Rows.Item(4).RowHeight = 53.4
$ws.Rows.Item(5).RowHeight = 40.2

# Update the active cell selection saved in the sheet view.
$ws.Range("D5").Select()
